$wb = $excel.ActiveWorkbook

# Update the translation text for "choose_authorization" on the
# table_specific_translations sheet: "Choose an Authorization" -> "Choose a Distribution"
$ws = $wb.Worksheets.Item("table_specific_translations")
$ws.Range("B3").Value = "Choose a Distribution"

# Make table_specific_translations the active sheet (was "properties")
$ws.Activate()

# Match the new active cell / selection on that sheet
$ws.Range("B4").Select()
